$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.446.51"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "3.276.55"
$ws.Range("E3").Value = "  +3.19%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "614.03"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").Value = "158.78"
$ws.Range("E6").Value = "  +2.98%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.275.81"
$ws.Range("E8").Value = "  +3.13%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  +3.39%  "

$ws.Range("E11").Value = "  +2.16%  "

$ws.Range("D12").Value = "0.498"
$ws.Range("E12").Value = "  -3.44%  "

$ws.Range("E13").Value = "  +3.23%  "

$ws.Range("D14").Value = "39.16"
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "3.814.12"
$ws.Range("E15").Value = "  +3.26%  "

$ws.Range("D16").Value = "66.501.82"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("D17").Value = "7.43"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("D18").Value = "3.272.54"
$ws.Range("E18").Value = "  +3.24%  "

$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "506.78"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "15.62"
$ws.Range("E21").Value = "  +2.11%  "

$ws.Range("D22").Value = "0.757"
$ws.Range("E22").Value = "  +4.10%  "

$ws.Range("E23").Value = "  +2.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.24%  "

$ws.Range("D25").Value = "86.78"
$ws.Range("E25").Value = "  +2.99%  "

$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  +2.35%  "

$ws.Range("D28").Value = "9.27"
$ws.Range("E28").Value = "  +1.38%  "

$ws.Range("E29").Value = "  +1.89%  "

$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +46.33%  "

$ws.Range("D31").Value = "7.02"
$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("E32").Value = "  -3.52%  "

$ws.Range("D33").Value = "28.17"
$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").Value = "6.51"
$ws.Range("E36").Value = "  +0.54%  "

$ws.Range("D37").Value = "3.53"
$ws.Range("E37").Value = "  +24.28%  "

$ws.Range("D38").Value = "0.0₃0794"
$ws.Range("E38").Value = "  +17.00%  "

$ws.Range("D39").Value = "55.52"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").Value = "499.51"
$ws.Range("E40").Value = "  -2.76%  "

$ws.Range("D41").Value = "0.0426"
$ws.Range("E41").Value = "  +1.86%  "

$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("E44").Value = "  +7.03%  "

$ws.Range("E45").Value = "  -0.32%  "

$ws.Range("D46").Value = "3.004.96"
$ws.Range("E46").Value = "  +6.65%  "

$ws.Range("D47").Value = "29.24"
$ws.Range("E47").Value = "  +5.08%  "

$ws.Range("D48").Value = "2.49"
$ws.Range("E48").Value = "  +5.56%  "

$ws.Range("E49").Value = "  +2.69%  "

$ws.Range("E50").Value = "  -0.01%  "

$ws.Range("D51").Value = "121.51"
$ws.Range("E51").Value = "  -0.16%  "
